# Update Fonds de solidarite dataset with 2022-06-14 data
# Only columns C (nombre_aides) and E (montant_total) change for the listed rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row=9;   C=286;    E=36710944}
    @{Row=24;  C=35711;  E=132418601}
    @{Row=37;  C=23051;  E=130209454}
    @{Row=91;  C=151205; E=482940711}
    @{Row=92;  C=409310; E=1597395258}
    @{Row=93;  C=209667; E=1310155630}
    @{Row=94;  C=94239;  E=919263386}
    @{Row=95;  C=50806;  E=934548578}
    @{Row=96;  C=17325;  E=797387398}
    @{Row=104; C=135303; E=272657154}
    @{Row=119; C=358;    E=10772714}
    @{Row=155; C=72711;  E=416714952}
    @{Row=165; C=83817;  E=355003773}
    @{Row=174; C=226111; E=900762198}
)

foreach ($u in $updates) {
    $ws.Range("C" + $u.Row).Value = $u.C
    $ws.Range("E" + $u.Row).Value = $u.E
}
